$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: shift the existing "Email" column (currently F) to H, for header + all 36 data rows ---
# (.Value doesn't read reliably in this COM layer, so use .Value2 for the read side.)
for ($r = 1; $r -le 37; $r++) {
    $src = $ws.Cells.Item($r, 6)
    $ws.Cells.Item($r, 8).Value = $src.Value2
}

# --- Step 2: new header cells F1 (submission date) and G1 (assignment name) ---
# Leading apostrophe forces text so "2020-11-09" is not auto-parsed into a date serial.
$ws.Cells.Item(1, 6).Value = "'2020-11-09"
$ws.Cells.Item(1, 7).Value = "Resenha Novos Clássicos"

# Re-apply the bold/bordered header look (lost on the newly-touched F1/G1/H1 cells) by
# copying the existing header formatting from A1 without disturbing the values just set.
$ws.Cells.Item(1, 1).Copy()
$ws.Cells.Item(1, 6).PasteSpecial(-4122)
$ws.Cells.Item(1, 1).Copy()
$ws.Cells.Item(1, 7).PasteSpecial(-4122)
$ws.Cells.Item(1, 1).Copy()
$ws.Cells.Item(1, 8).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Step 3: grade data for the new "Resenha Novos Clássicos" assignment ---
# Column F = score value, column G = grade. Rows with no entry are left blank (matches source).
    $ws.Cells.Item(2, 6).Value = 19.47
    $ws.Cells.Item(2, 7).Value = 5.0
    $ws.Cells.Item(3, 6).Value = 6.19
    $ws.Cells.Item(3, 7).Value = 5.0
    $ws.Cells.Item(4, 6).Value = 93.81
    $ws.Cells.Item(4, 7).Value = 7.0
    $ws.Cells.Item(5, 6).Value = 54.87
    $ws.Cells.Item(5, 7).Value = 5.0
    $ws.Cells.Item(6, 6).Value = 91.15
    $ws.Cells.Item(6, 7).Value = 7.0
    $ws.Cells.Item(7, 6).Value = 0.0
    $ws.Cells.Item(7, 7).Value = 0.0
    $ws.Cells.Item(8, 6).Value = 77.88
    $ws.Cells.Item(9, 6).Value = 92.92
    $ws.Cells.Item(9, 7).Value = 5.0
    $ws.Cells.Item(10, 6).Value = 76.11
    $ws.Cells.Item(10, 7).Value = 0.0
    $ws.Cells.Item(11, 6).Value = 51.33
    $ws.Cells.Item(11, 7).Value = 5.0
    $ws.Cells.Item(12, 6).Value = 0.0
    $ws.Cells.Item(12, 7).Value = 0.0
    $ws.Cells.Item(13, 6).Value = 0.0
    $ws.Cells.Item(13, 7).Value = 0.0
    $ws.Cells.Item(14, 6).Value = 74.34
    $ws.Cells.Item(14, 7).Value = 5.0
    $ws.Cells.Item(15, 7).Value = 0.0
    $ws.Cells.Item(16, 6).Value = 16.81
    $ws.Cells.Item(16, 7).Value = 0.0
    $ws.Cells.Item(17, 7).Value = 3.0
    $ws.Cells.Item(18, 6).Value = 91.15
    $ws.Cells.Item(18, 7).Value = 5.0
    $ws.Cells.Item(19, 6).Value = 82.3
    $ws.Cells.Item(19, 7).Value = 7.0
    $ws.Cells.Item(20, 6).Value = 93.81
    $ws.Cells.Item(20, 7).Value = 5.0
    $ws.Cells.Item(21, 6).Value = 14.16
    $ws.Cells.Item(21, 7).Value = 7.0
    $ws.Cells.Item(22, 6).Value = 91.15
    $ws.Cells.Item(22, 7).Value = 5.0
    $ws.Cells.Item(23, 6).Value = 18.58
    $ws.Cells.Item(23, 7).Value = 10.0
    $ws.Cells.Item(24, 6).Value = 93.81
    $ws.Cells.Item(24, 7).Value = 10.0
    $ws.Cells.Item(25, 6).Value = 0.0
    $ws.Cells.Item(26, 6).Value = 0.88
    $ws.Cells.Item(26, 7).Value = 7.0
    $ws.Cells.Item(27, 6).Value = 0.0
    $ws.Cells.Item(27, 7).Value = 0.0
    $ws.Cells.Item(29, 7).Value = 0.0
    $ws.Cells.Item(30, 6).Value = 4.42
    $ws.Cells.Item(30, 7).Value = 7.0
    $ws.Cells.Item(31, 6).Value = 0.0
    $ws.Cells.Item(31, 7).Value = 7.0
    $ws.Cells.Item(32, 6).Value = 61.06
    $ws.Cells.Item(32, 7).Value = 7.0
    $ws.Cells.Item(33, 6).Value = 30.09
    $ws.Cells.Item(33, 7).Value = 7.0
    $ws.Cells.Item(34, 6).Value = 0.0
    $ws.Cells.Item(34, 7).Value = 5.0
    $ws.Cells.Item(35, 6).Value = 70.8
    $ws.Cells.Item(35, 7).Value = 7.0
    $ws.Cells.Item(36, 6).Value = 2.65
    $ws.Cells.Item(36, 7).Value = 0.0
    $ws.Cells.Item(37, 7).Value = 0.0
